$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (the 0.5W Chip LED line) gets a new manufacturer / supplier part number.
# Leading "'" preserves the existing quote-prefixed cell style (s=2) instead of
# letting Excel drop it to the no-quote-prefix style when the text is rewritten.
$ws.Range("H4").Value = "'365-1546-1-ND"

# D4 (manufacturer part number) gets new text AND a smaller, unbordered font.
$ws.Range("D4").Value = "OVS5MRBCR4"
$ws.Range("D4").Style = "Normal"
$ws.Range("D4").Font.Size = 9

$ws.Range("C4").Value = "'TT Electronics/Optek Technology"

# Updated unit price and its extended (qty * unit price) subtotal.
$ws.Range("I4").Value = 0.63
$ws.Range("J4").Value = 12.6

# Leave the cursor on A2, matching the saved selection in the workbook.
$ws.Range("A2").Select() | Out-Null
